# Auto-generated edit script: updates static market-price derived cells
# across all 8 sheets, per the scheduled-runner data refresh diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1599.5
$ws.Range("I4").Value = 1199
$ws.Range("K4").Value = 1199
$ws.Range("M4").Value = -1085
$ws.Range("H62").Value = 5533
$ws.Range("I62").Value = 5200
$ws.Range("K62").Value = 5200
$ws.Range("M62").Value = -4576
$ws.Range("H64").Value = 5272.0376
$ws.Range("I64").Value = 5119.886
$ws.Range("K64").Value = 5119.886
$ws.Range("M64").Value = -4871.886
$ws.Range("H65").Value = 5533
$ws.Range("I65").Value = 5200
$ws.Range("K65").Value = 26000
$ws.Range("M65").Value = -22880
$ws.Range("H67").Value = 5272.0376
$ws.Range("I67").Value = 5119.886
$ws.Range("K67").Value = 5119.886
$ws.Range("M67").Value = -4261.886
$ws.Range("H80").Value = 1479
$ws.Range("I80").Value = 401
$ws.Range("J80").Value = 1969
$ws.Range("K80").Value = 1203
$ws.Range("L80").Value = 5907
$ws.Range("M80").Value = -205
$ws.Range("N80").Value = -7903
$ws.Range("H83").Value = 1479
$ws.Range("I83").Value = 401
$ws.Range("J83").Value = 1969
$ws.Range("K83").Value = 3609
$ws.Range("L83").Value = 17721
$ws.Range("M83").Value = 1383
$ws.Range("N83").Value = -27705
$ws.Range("H98").Value = 280
$ws.Range("I98").Value = 280
$ws.Range("K98").Value = 280
$ws.Range("M98").Value = 1218
$ws.Range("H122").Value = 280
$ws.Range("I122").Value = 280
$ws.Range("K122").Value = 840
$ws.Range("M122").Value = 1610

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9408.548000000001
$ws.Range("I32").Value = 3572.4167
$ws.Range("J32").Value = 17190.055
$ws.Range("K32").Value = 3572.4167
$ws.Range("L32").Value = 17190.055
$ws.Range("M32").Value = -3285.4167
$ws.Range("N32").Value = -17764.055
$ws.Range("H61").Value = 85848.836
$ws.Range("I61").Value = 2744.2727
$ws.Range("K61").Value = 2744.2727
$ws.Range("M61").Value = -2532.2727
$ws.Range("H132").Value = 1829.5238
$ws.Range("I132").Value = 1682.2059
$ws.Range("K132").Value = 5046.6177
$ws.Range("M132").Value = -2516.6177
$ws.Range("H136").Value = 85848.836
$ws.Range("I136").Value = 2744.2727
$ws.Range("K136").Value = 8232.8181
$ws.Range("M136").Value = -5682.8181

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7694520.5
$ws.Range("I107").Value = 10002010
$ws.Range("K107").Value = 10002010
$ws.Range("M107").Value = -10000090
$ws.Range("H132").Value = 31366
$ws.Range("J132").Value = 31366
$ws.Range("L132").Value = 31366
$ws.Range("N132").Value = -41486

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1549.75
$ws.Range("I12").Value = 1999.6666
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 1999.6666
$ws.Range("L12").Value = 200
$ws.Range("M12").Value = -1829.6666
$ws.Range("N12").Value = -540
$ws.Range("H13").Value = 121.5
$ws.Range("J13").Value = 121.5
$ws.Range("L13").Value = 121.5
$ws.Range("N13").Value = -399.5
$ws.Range("H58").Value = 1645
$ws.Range("I58").Value = 1466.1333
$ws.Range("K58").Value = 1466.1333
$ws.Range("M58").Value = -1263.1333
$ws.Range("H93").Value = 27120
$ws.Range("I93").Value = 27120
$ws.Range("K93").Value = 27120
$ws.Range("M93").Value = -25248
$ws.Range("H107").Value = 1802.826
$ws.Range("I107").Value = 1511.6923
$ws.Range("J107").Value = 2181.3
$ws.Range("K107").Value = 1511.6923
$ws.Range("L107").Value = 2181.3
$ws.Range("M107").Value = 408.3077000000001
$ws.Range("N107").Value = -6021.3
$ws.Range("H134").Value = 4026531.5
$ws.Range("I134").Value = 4467145
$ws.Range("K134").Value = 13401435
$ws.Range("M134").Value = -13398900
$ws.Range("H136").Value = 1645
$ws.Range("I136").Value = 1466.1333
$ws.Range("K136").Value = 4398.3999
$ws.Range("M136").Value = -1848.3999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 10000
$ws.Range("I82").Value = 10000
$ws.Range("K82").Value = 30000
$ws.Range("M82").Value = -29594
$ws.Range("H85").Value = 10000
$ws.Range("I85").Value = 10000
$ws.Range("K85").Value = 30000
$ws.Range("M85").Value = -28596

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 7789453.5
$ws.Range("I11").Value = 2186409
$ws.Range("J11").Value = 13392498
$ws.Range("K11").Value = 2186409
$ws.Range("L11").Value = 13392498
$ws.Range("M11").Value = -2186270
$ws.Range("N11").Value = -13392776
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H80").Value = 1549.625
$ws.Range("I80").Value = 999.25
$ws.Range("J80").Value = 2100
$ws.Range("K80").Value = 999.25
$ws.Range("L80").Value = 2100
$ws.Range("M80").Value = -1.25
$ws.Range("N80").Value = -4096
$ws.Range("H83").Value = 1549.625
$ws.Range("I83").Value = 999.25
$ws.Range("J83").Value = 2100
$ws.Range("K83").Value = 4996.25
$ws.Range("L83").Value = 10500
$ws.Range("M83").Value = -4.25
$ws.Range("N83").Value = -20484
$ws.Range("H93").Value = 17105.666
$ws.Range("J93").Value = 17105.666
$ws.Range("L93").Value = 17105.666
$ws.Range("N93").Value = -20849.666
$ws.Range("H107").Value = 571.1
$ws.Range("I107").Value = 555.9048
$ws.Range("K107").Value = 555.9048
$ws.Range("M107").Value = 1364.0952
$ws.Range("H122").Value = 8635670
$ws.Range("I122").Value = 14030627
$ws.Range("J122").Value = 3737.8
$ws.Range("K122").Value = 42091881
$ws.Range("L122").Value = 11213.4
$ws.Range("M122").Value = -42089431
$ws.Range("N122").Value = -16113.4
$ws.Range("H132").Value = 3893.88
$ws.Range("I132").Value = 3064.1428
$ws.Range("J132").Value = 8250
$ws.Range("K132").Value = 9192.428400000001
$ws.Range("L132").Value = 24750
$ws.Range("M132").Value = -6662.428400000001
$ws.Range("N132").Value = -29810

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1170.8695
$ws.Range("I22").Value = 1329.3077
$ws.Range("J22").Value = 964.9
$ws.Range("K22").Value = 1329.3077
$ws.Range("L22").Value = 964.9
$ws.Range("M22").Value = -1034.3077
$ws.Range("N22").Value = -1554.9
$ws.Range("H27").Value = 1170.8695
$ws.Range("I27").Value = 1329.3077
$ws.Range("J27").Value = 964.9
$ws.Range("K27").Value = 1329.3077
$ws.Range("L27").Value = 964.9
$ws.Range("M27").Value = -1222.3077
$ws.Range("N27").Value = -1178.9
$ws.Range("H34").Value = 1133
$ws.Range("I34").Value = 1133
$ws.Range("K34").Value = 1133
$ws.Range("M34").Value = -961
$ws.Range("H40").Value = 11115539
$ws.Range("I40").Value = 4583.1665
$ws.Range("J40").Value = 19448756
$ws.Range("K40").Value = 4583.1665
$ws.Range("L40").Value = 19448756
$ws.Range("M40").Value = -4447.1665
$ws.Range("N40").Value = -19449028
$ws.Range("H61").Value = 5624
$ws.Range("I61").Value = 4999
$ws.Range("K61").Value = 4999
$ws.Range("M61").Value = -4797
$ws.Range("H68").Value = 192763
$ws.Range("I68").Value = 235332.56
$ws.Range("J68").Value = 1200
$ws.Range("K68").Value = 235332.56
$ws.Range("L68").Value = 1200
$ws.Range("M68").Value = -234583.56
$ws.Range("N68").Value = -2698
$ws.Range("H71").Value = 192763
$ws.Range("I71").Value = 235332.56
$ws.Range("J71").Value = 1200
$ws.Range("K71").Value = 1176662.8
$ws.Range("L71").Value = 6000
$ws.Range("M71").Value = -1172918.8
$ws.Range("N71").Value = -13488
$ws.Range("H82").Value = 1998.8182
$ws.Range("I82").Value = 1927.4286
$ws.Range("J82").Value = 2123.75
$ws.Range("K82").Value = 1927.4286
$ws.Range("L82").Value = 2123.75
$ws.Range("M82").Value = -1566.4286
$ws.Range("N82").Value = -2845.75
$ws.Range("H85").Value = 1998.8182
$ws.Range("I85").Value = 1927.4286
$ws.Range("J85").Value = 2123.75
$ws.Range("K85").Value = 1927.4286
$ws.Range("L85").Value = 2123.75
$ws.Range("M85").Value = -679.4286
$ws.Range("N85").Value = -4619.75
$ws.Range("H113").Value = 5624
$ws.Range("I113").Value = 4999
$ws.Range("K113").Value = 4999
$ws.Range("M113").Value = -2829
$ws.Range("H122").Value = 65003516
$ws.Range("I122").Value = 62503620
$ws.Range("K122").Value = 187510860
$ws.Range("M122").Value = -187508410
$ws.Range("H132").Value = 2254
$ws.Range("I132").Value = 2254
$ws.Range("K132").Value = 6762
$ws.Range("M132").Value = -4232

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 99450
$ws.Range("I57").Value = 78900
$ws.Range("K57").Value = 78900
$ws.Range("M57").Value = -78146
$ws.Range("H81").Value = 1763
$ws.Range("J81").Value = 1763
$ws.Range("L81").Value = 3526
$ws.Range("N81").Value = -5648
$ws.Range("H84").Value = 1763
$ws.Range("J84").Value = 1763
$ws.Range("L84").Value = 17630
$ws.Range("N84").Value = -28238
$ws.Range("H107").Value = 3620.4146
$ws.Range("J107").Value = 2212.75
$ws.Range("L107").Value = 6638.25
$ws.Range("N107").Value = -10478.25
$ws.Range("H122").Value = 1791.6086
$ws.Range("I122").Value = 1193.9333
$ws.Range("K122").Value = 3581.7999
$ws.Range("M122").Value = -1131.7999
$ws.Range("H132").Value = 2605.4583
$ws.Range("I132").Value = 2357.2778
$ws.Range("K132").Value = 7071.8334
$ws.Range("M132").Value = -4541.8334
$ws.Range("H136").Value = 1704.4615
$ws.Range("I136").Value = 1290.8
$ws.Range("K136").Value = 3872.4
$ws.Range("M136").Value = -1322.4
